# Support Vector Machine training implemented
# Add a new "CompleteSVM" worksheet after the existing "CompleteRFC" sheet,
# populate it with the SVM configuration (Kernel/Gamma/C = rbf/scale/1.9),
# and make it the active sheet (matching tabSelected/activeTab move in the
# workbook XML).

$wb = $excel.ActiveWorkbook
$wsRFC = $wb.Worksheets.Item("CompleteRFC")

# Insert the new worksheet right after CompleteRFC.
$wsSVM = $wb.Worksheets.Add($null, $wsRFC)
$wsSVM.Name = "CompleteSVM"

# Header row
$wsSVM.Range("A1").Value = "Kernel"
$wsSVM.Range("B1").Value = "Gamma"
$wsSVM.Range("C1").Value = "C"

# Data row
$wsSVM.Range("A2").Value = "rbf"
$wsSVM.Range("B2").Value = "scale"
$wsSVM.Range("C2").Value = 1.9

# Match the selection left on the new sheet in the original edit.
$wsSVM.Range("C4").Select() | Out-Null

# CompleteSVM becomes the active/visible tab.
$wsSVM.Activate() | Out-Null
